$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(6, 15, 5, 5),
    @(2, 16, 1, 4),
    @(4, 2, 8, 18),
    @(7, 15, 6, 5),
    @(4, 4, 5, 16),
    @(3, 6, 4, 14),
    @(5, 15, 4, 5),
    @(5, 13, 4, 7),
    @(2, 7, 5, 13),
    @(5, 15, 7, 5),
    @(4, 12, 3, 8),
    @(5, 7, 3, 13),
    @(4, 3, 6, 17),
    @(3, 12, 2, 8),
    @(5, 13, 6, 7),
    @(8, 5, 7, 15),
    @(4, 12, 5, 8),
    @(2, 6, 3, 14),
    @(2, 4, 3, 16),
    @(7, 5, 6, 15),
    @(5, 4, 7, 16)
)

$startRow = 1554
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
}

$ws.Range("A1575").Select()
$excel.ActiveWindow.ScrollRow = 1559
$excel.ActiveWindow.ScrollColumn = 1

